$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for rows 3, 4 and 5 in the columns
# that change (D, J, K, L, M, P). The edit performs a cyclic rotation of
# these values: row3 -> row4, row4 -> row5, row5 -> row3.
$cols = @("D", "J", "K", "L", "M", "P")

$row3 = @{}
$row4 = @{}
$row5 = @{}

foreach ($col in $cols) {
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
    $row5[$col] = $ws.Range("${col}5").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}4").Value2 = $row3[$col]
    $ws.Range("${col}5").Value2 = $row4[$col]
    $ws.Range("${col}3").Value2 = $row5[$col]
}
